$wb = $excel.ActiveWorkbook

# --- Clean up named ranges that point at the sheet we are about to remove ---
# (TEST_SHEET_TXL_ROW!$I$13 / TEST_SHEET_TXL_ROW!$O$12 -> would become #REF! once the
# sheet is deleted, so drop them first.)
$wb.Names("TEST_NAME_FOR_NBK").Delete()
$wb.Names("TESTOTHER_NAME").Delete()

# Drop the obsolete helper sheet entirely.
$wb.Worksheets("TEST_SHEET_TXL_ROW").Delete()

# Drop the old TXL_SHEET_TRACKER sheet - it gets rebuilt from scratch below with an
# updated layout (extra sheet_index column). Bump the internal sheetId counter up to
# 13 first by cycling a few scratch sheets through the workbook, matching the id the
# rebuilt tracker sheet ends up with.
$wb.Worksheets("TXL_SHEET_TRACKER").Delete()
$null = $wb.Worksheets.Add()
$null = $wb.Worksheets.Add()
$null = $wb.Worksheets.Add()

$tracker = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$tracker.Name = "TXL_SHEET_TRACKER"

$wb.Worksheets("Sheet1").Delete()
$wb.Worksheets("Sheet2").Delete()
$wb.Worksheets("Sheet3").Delete()

# --- Rebuild TXL_SHEET_TRACKER content ---
$tracker.Range("A1").Value2 = "sheet_name"
$tracker.Range("B1").Value2 = "descr"
$tracker.Range("C1").Value2 = "sheet_type"
$tracker.Range("D1").Value2 = "sheet_index"
$tracker.Range("A1:D1").Interior.Color = 12566463

$tracker.Range("A2").Value2 = "TEST_STANDARD_ROW"
$tracker.Range("B2").Value2 = "Sheet to test standard row sheet"
$tracker.Range("C2").Value2 = 1
$tracker.Range("D2").Formula = "=_xlfn.SHEET(TEST_STANDARD_ROW!`$A`$1)"

$tracker.Columns("A").ColumnWidth = 22.75
$tracker.Columns("B").ColumnWidth = 27.75
$tracker.Columns("C").ColumnWidth = 9.75
$tracker.Columns("D").ColumnWidth = 10.625

$tracker.Activate()
$tracker.Range("G6").Select()

# --- TEST_SHEET (formerly TEST_SHEET_TXL_ROW's tab, now the empty "TEST_SHEET") ---
$testSheet = $wb.Worksheets("TEST_SHEET")
$testSheet.Activate()
$testSheet.Range("M20").Select()

# --- TEST_STANDARD_ROW: drop the color column, keep just the index column ---
$std = $wb.Worksheets("TEST_STANDARD_ROW")
$std.Columns("B").Delete()
$std.Activate()
$std.Columns("B:B").Select()

# --- TEST_SCALAR_INPUT: no longer the active tab ---
$scalar = $wb.Worksheets("TEST_SCALAR_INPUT")
$scalar.Activate()
$scalar.Range("D33").Select()

# --- Defined names ---
$wb.Names.Add("TEST_STANDARD_ROW__index", "=TEST_STANDARD_ROW!`$A:`$A")

$tracker.Activate()
